$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A202").NumberFormat = "@"
$ws.Range("A202").Value = "01-09-2021"
$ws.Range("A202").ClearFormats()
$ws.Range("B202:I202").Value = 0
$ws.Range("J202").Value = 1891
$ws.Range("K202:L202").Value = 0
$ws.Range("M202").Value = 1891
$ws.Range("N202:P202").Value = 0
$ws.Range("Q202").Value = -1891
$ws.Range("R202:V202").Value = 0
$ws.Range("W202").Value = 1891
$ws.Range("X202").Value = -1891

$ws.Range("J194").Value = 2381
$ws.Range("M194").Value = 2381
$ws.Range("Q194").Value = -2381
$ws.Range("W194").Value = 2381
$ws.Range("X194").Value = -2381

$ws.Range("J195").Value = 2338
$ws.Range("M195").Value = 2338
$ws.Range("Q195").Value = -2338
$ws.Range("W195").Value = 2338
$ws.Range("X195").Value = -2338

$ws.Range("J196").Value = 2294
$ws.Range("M196").Value = 2294
$ws.Range("Q196").Value = -2294
$ws.Range("W196").Value = 2294
$ws.Range("X196").Value = -2294

$ws.Range("J197").Value = 2211
$ws.Range("M197").Value = 2211
$ws.Range("Q197").Value = -2211
$ws.Range("W197").Value = 2211
$ws.Range("X197").Value = -2211

$ws.Range("J198").Value = 2148
$ws.Range("M198").Value = 2148
$ws.Range("Q198").Value = -2148
$ws.Range("W198").Value = 2148
$ws.Range("X198").Value = -2148

$ws.Range("J199").Value = 2085
$ws.Range("M199").Value = 2085
$ws.Range("Q199").Value = -2085
$ws.Range("W199").Value = 2085
$ws.Range("X199").Value = -2085

$ws.Range("J200").Value = 2020
$ws.Range("M200").Value = 2020
$ws.Range("Q200").Value = -2020
$ws.Range("W200").Value = 2020
$ws.Range("X200").Value = -2020

$ws.Range("J201").Value = 1955
$ws.Range("M201").Value = 1955
$ws.Range("Q201").Value = -1955
$ws.Range("W201").Value = 1955
$ws.Range("X201").Value = -1955
